$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# D12: numeric value 100% formatted as percentage
$ws.Range("D12").Value = 1
$ws.Range("D12").NumberFormat = "0%"

# E12: new progress note
$ws.Range("E12").Value = "100%(17/06/2010)"

# Update the selection / view state
$ws.Range("G14").Select()
